$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function TimeFrac([int]$h, [int]$m) {
    return ($h + $m/60)/24
}

# Update the time-of-day values in D3:D14. (D2 is left as-is; D10 used to
# hold a stray 18:07 value that's now back in the ascending 8:40->9:40
# sequence the rest of the column follows.)
$ws.Range("D3").Value  = TimeFrac 8 40
$ws.Range("D4").Value  = TimeFrac 9 2
$ws.Range("D5").Value  = TimeFrac 9 5
$ws.Range("D6").Value  = TimeFrac 9 11
$ws.Range("D7").Value  = TimeFrac 9 16
$ws.Range("D8").Value  = TimeFrac 9 20
$ws.Range("D9").Value  = TimeFrac 9 22
$ws.Range("D10").Value = TimeFrac 9 25
$ws.Range("D11").Value = TimeFrac 9 28
$ws.Range("D12").Value = TimeFrac 9 32
$ws.Range("D13").Value = TimeFrac 9 36
$ws.Range("D14").Value = TimeFrac 9 40

# Move the active selection from D11 down to D15 (this also drops the
# stale topLeftCell="A8" scroll anchor left over from the previous view).
$ws.Range("D15").Select()
